$d = $word.ActiveDocument

# Locate the paragraph that currently ends with (tab +) "The " -- this is
# the paragraph right after "Description:" where the new description text
# needs to be appended.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "The \r$") {
        $target = $p
        break
    }
}

$r = $target.Range
# Collapse to an insertion point right before the paragraph mark, i.e.
# immediately after "The ".
$ip = $d.Range($r.End - 1, $r.End - 1)

$ip.InsertAfter("Server manager, or Server Server, is an application that can be used to h")
$ip = $d.Range($ip.End, $ip.End)

$ip.InsertAfter("o")
$ip = $d.Range($ip.End, $ip.End)

$ip.InsertAfter("s")
$ip = $d.Range($ip.End, $ip.End)

$ip.InsertAfter("t")
$ip = $d.Range($ip.End, $ip.End)

$ip.InsertAfter(" multiple game servers at once and allow control of all of them remotely though a single application instead of multiple.")
$ip = $d.Range($ip.End, $ip.End)

$ip.InsertAfter(" The Server Server is controlled though issuing commands through a console. In the console info about the different game servers is visible and their outputs. Multiple users can interact with the Server Server simultaneously and only have certain permissions that limit their ability to use commands. ")
$ip = $d.Range($ip.End, $ip.End)

$ip.InsertAfter("Commands can also be created/edited in real time without the entire application needing to restart.")
$ip = $d.Range($ip.End, $ip.End)

# Add the new, completely blank paragraph that separates the description
# from the "Gantt Chart:" paragraph. InsertXML with a bare <w:p/> yields a
# true empty paragraph (no placeholder run), matching a plain
# Enter-keypress split.
$blankParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
[void]$ip.InsertXML($blankParaXml)
